$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D price updates (row => new value)
$dUpdates = @{
    2 = "281.39"
    3 = "20.92"
    4 = "6.258"
    5 = "0.06144"
    8 = "1.497"
    9 = "0.8171"
    10 = "0.01382"
    12 = "0.08330"
    14 = "0.03173"
    15 = "0.09135"
    16 = "3.709"
    17 = "0.001643"
    18 = "0.04676"
    19 = "0.006454"
    23 = "3.796"
    24 = "2.337"
    25 = "0.3367"
    40 = "0.04668"
    41 = "0.005603"
    42 = "0.007163"
    43 = "0.1098"
    44 = "0.01111"
    45 = "0.00006154"
    46 = "0.00000000750"
    48 = "0.002947"
}

foreach ($row in $dUpdates.Keys) {
    $cell = $ws.Range("D" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $dUpdates[$row]
}

# Column G (Hora) updates: every data row (2-51) goes from "7" to "8"
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Range("G" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = "8"
}
